$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 values ---
$ws.Cells.Item(2, 2).Value = "java"
$ws.Cells.Item(2, 4).Value = "Juan Perez"
$ws.Cells.Item(2, 5).Value = 1450000.0
$ws.Cells.Item(2, 6).Value = "Carlos Muñoz"
$ws.Cells.Item(2, 7).Value = 5.0

# --- Update row 3 values ---
$ws.Cells.Item(3, 1).Value = 2222.0
$ws.Cells.Item(3, 3).Value = 70.0
$ws.Cells.Item(3, 4).Value = "Jorge mundaca"
$ws.Cells.Item(3, 5).Value = 1450000.0
$ws.Cells.Item(3, 6).Value = "Claudia nog"
$ws.Cells.Item(3, 7).Value = 6.666666666666667
$ws.Cells.Item(3, 8).Value = "SF: AA"

# --- Remove row 4 entirely (table now only spans rows 1-3) ---
$ws.Range("A4:H4").Delete()

# --- Adjust column widths for D and F to match new (longer) content ---
$ws.Columns.Item(4).ColumnWidth = 13.666666666666666
$ws.Columns.Item(6).ColumnWidth = 12.499999999999998
